$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.200496077537537
$ws.Range("B1").Value = 2.061381101608276
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.080751657485962
$ws.Range("E1").Value = 1.206899642944336
